$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.761.28'
$ws.Range("E2").Value = '  +4.01%  '
$ws.Range("D3").Value = '2.262.40'
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.531'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.42%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.481'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.69'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0794'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("D15").Value = '2.615.04'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("D17").Value = '2.282.99'
$ws.Range("E17").Value = '  +2.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.760'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.22%  '
$ws.Range("D19").Value = '41.659.62'
$ws.Range("E19").Value = '  +3.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.86%  '
$ws.Range("D21").Value = '0.0₃0903'
$ws.Range("E21").Value = '  +1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.08%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +4.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("E29").Value = '  +11.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.02%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0743'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.54%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.76'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.20%  '
$ws.Range("E39").Value = '  +2.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.103'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.51%  '
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.57%  '
$ws.Range("D43").Value = '2.059.97'
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0277'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.88'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.32%  '
$ws.Range("E48").Value = '  +6.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.87%  '
